$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp footer in A1
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 16:22"

# Full country table (rows 4-215), re-sorted by "Casos totales" (column B) descending
# after applying the latest per-country case updates.
$data = @(
    ,@('Estados Unidos',400549,214,21711,365981,9169,16,12857)
    ,@('España',146690,4748,48021,84114,7069,510,14555)
    ,@('Italia',135586,0,24392,94067,3792,0,17127)
    ,@('Alemania',109329,1666,36081,71152,4895,80,2096)
    ,@('Francia',109069,0,19337,79404,7131,0,10328)
    ,@('China',81802,62,77279,1190,189,2,3333)
    ,@('Iran',64586,1997,27039,33554,3987,121,3993)
    ,@('Reino Unido',55242,0,135,48948,1559,0,6159)
    ,@('Turquia',34109,0,1582,31802,1474,0,725)
    ,@('Belgica',23403,1209,4681,16482,1276,205,2240)
    ,@('Suiza',22789,536,8704,13227,391,37,858)
    ,@('Paises Bajos',20549,969,250,18051,1424,147,2248)
    ,@('Canada',17897,0,4050,13466,426,0,381)
    ,@('Brasil',14152,118,127,13326,296,13,699)
    ,@('Portugal',13141,699,196,12565,245,35,380)
    ,@('Austria',12881,242,4512,8096,267,30,273)
    ,@('Corea del Sur',10384,53,6776,3408,55,8,200)
    ,@('Israel',9404,156,801,8531,147,7,72)
    ,@('Rusia',8672,1175,580,8029,8,5,63)
    ,@('Suecia',8419,726,205,7527,678,96,687)
    ,@('Noruega',6086,0,32,5961,78,4,93)
    ,@('Australia',6013,25,2813,3150,87,1,50)
    ,@('Irlanda',5709,0,25,5474,165,0,210)
    ,@('India',5480,129,468,4848,0,4,164)
    ,@('Dinamarca',5386,315,1621,3547,127,15,218)
    ,@('Chile',5116,0,898,4175,337,0,43)
    ,@('Chequia',5033,16,181,4761,103,3,91)
    ,@('Polonia',5000,152,222,4642,160,7,136)
    ,@('Rumania',4761,344,528,4018,162,18,215)
    ,@('Japon',4257,0,622,3542,80,0,93)
    ,@('Pakistan',4183,148,467,3658,25,1,58)
    ,@('Malasia',4119,156,1487,2567,76,2,65)
    ,@('Ecuador',3995,0,140,3635,156,0,220)
    ,@('Filipinas',3870,106,96,3592,1,5,182)
    ,@('Luxemburgo',2970,0,500,2426,35,0,44)
    ,@('Indonesia',2956,218,222,2494,0,19,240)
    ,@('Peru',2954,0,1301,1546,109,0,107)
    ,@('Arabia Saudita',2932,137,615,2276,41,0,41)
    ,@('Mexico',2785,346,633,2011,89,16,141)
    ,@('Serbia',2666,219,118,2483,112,4,65)
    ,@('Emiratos Arabes Unidos',2659,300,239,2408,1,0,12)
    ,@('Finlandia',2487,179,300,2147,82,6,40)
    ,@('Tailandia',2369,111,888,1451,61,3,30)
    ,@('Panama',2249,149,16,2174,91,4,59)
    ,@('Catar',2210,153,178,2026,37,0,6)
    ,@('Republica Dominicana',1956,0,36,1822,147,0,98)
    ,@('Grecia',1832,0,269,1482,90,0,81)
    ,@('Colombia',1780,0,100,1630,76,0,50)
    ,@('Sudafrica',1749,0,95,1641,7,0,13)
    ,@('Argentina',1715,0,338,1315,96,2,62)
    ,@('Ucrania',1668,206,35,1581,33,7,52)
    ,@('Singapur',1623,142,406,1211,29,0,6)
    ,@('Islandia',1616,30,633,977,13,0,6)
    ,@('Argelia',1468,0,113,1162,46,0,193)
    ,@('Egipto',1450,0,276,1080,0,0,94)
    ,@('Croacia',1343,61,179,1145,36,1,19)
    ,@('Marruecos',1242,58,97,1054,1,1,91)
    ,@('Nueva Zelanda',1210,50,282,927,4,0,1)
    ,@('Estonia',1185,36,72,1089,11,3,24)
    ,@('Irak',1122,0,373,684,0,0,65)
    ,@('Eslovenia',1091,32,120,931,35,4,40)
    ,@('Bielorrusia',1066,205,77,976,33,0,13)
    ,@('Moldavia',1056,0,40,992,80,2,24)
    ,@('Hong Kong',961,25,264,693,14,0,4)
    ,@('Lituania',912,32,8,889,21,0,15)
    ,@('Hungria',895,78,94,743,17,11,58)
    ,@('Armenia',881,28,114,758,30,1,9)
    ,@('Kuwait',855,112,111,743,21,0,1)
    ,@('Barein',821,10,467,349,3,0,5)
    ,@('Bosnia y Herzegovina',799,35,79,686,4,1,34)
    ,@('Azerbaiyan',717,0,44,665,23,0,8)
    ,@('Crucero',712,0,619,82,10,0,11)
    ,@('Kazajistan',709,12,54,648,21,1,7)
    ,@('Camerun',685,0,60,616,0,0,9)
    ,@('Eslovaquia',682,101,13,667,3,0,2)
    ,@('Tunez',623,0,25,575,67,0,23)
    ,@('Republica de Macedonia',617,18,35,553,15,3,29)
    ,@('Bulgaria',581,4,42,516,21,0,23)
    ,@('Letonia',577,29,16,559,5,0,2)
    ,@('Libano',575,27,62,494,29,0,19)
    ,@('Principado de Andorra',564,19,52,489,17,1,23)
    ,@('Uzbekistan',545,25,30,512,8,1,3)
    ,@('Republica de Chipre',494,0,47,438,13,0,9)
    ,@('Costa Rica',483,0,24,457,14,0,2)
    ,@('Uruguay',424,0,150,267,14,0,7)
    ,@('Afganistan',423,0,18,391,0,0,14)
    ,@('Oman',419,48,72,345,3,0,2)
    ,@('Albania',400,17,154,224,7,0,22)
    ,@('Cuba',396,0,27,358,15,0,11)
    ,@('Burkina Faso',384,0,127,238,0,0,19)
    ,@('Taiwan',379,3,67,307,0,0,5)
    ,@('Reunion',358,0,40,318,4,0,0)
    ,@('Jordania',353,0,138,209,5,0,6)
    ,@('Costa de Marfil',349,0,41,305,0,0,3)
    ,@('Honduras',312,7,6,284,10,0,22)
    ,@('Malta',299,6,5,294,4,0,0)
    ,@('Ghana',287,0,31,251,2,0,5)
    ,@('San Marino',279,0,40,205,14,0,34)
    ,@('Niger',278,0,26,241,0,0,11)
    ,@('Kirguistan',270,42,33,233,5,0,4)
    ,@('Mauricio',268,0,8,253,3,0,7)
    ,@('Estado de Palestina',263,2,44,218,0,0,1)
    ,@('Nigeria',254,0,44,204,2,0,6)
    ,@('Vietnam',251,0,126,125,8,0,0)
    ,@('Montenegro',248,7,4,242,7,0,2)
    ,@('Senegal',244,7,113,129,1,0,2)
    ,@('Banglades',218,54,33,165,1,3,20)
    ,@('Bolivia',210,16,2,193,3,1,15)
    ,@('Georgia',208,12,48,157,6,0,3)
    ,@('Sri Lanka',189,4,44,138,5,1,7)
    ,@('Islas Feroe',184,0,131,53,1,0,0)
    ,@('Consejo Danes para los Refugiados',180,0,9,153,0,0,18)
    ,@('Kenia',179,7,9,164,2,0,6)
    ,@('Mayotte',171,0,22,147,3,0,2)
    ,@('Venezuela',166,0,65,94,6,0,7)
    ,@('Isla de Man',158,8,80,77,7,0,1)
    ,@('Martinica',152,0,50,98,20,0,4)
    ,@('Guinea',144,0,5,139,0,0,0)
    ,@('Guadalupe',139,0,31,101,14,0,7)
    ,@('Republica de Yibuti',135,45,25,110,0,0,0)
    ,@('Brunei',135,0,91,43,3,0,1)
    ,@('Paraguay',119,4,15,99,1,0,5)
    ,@('Camboya',117,2,63,54,1,0,0)
    ,@('Gibraltar',113,0,60,53,0,0,0)
    ,@('Trinidad yTobago',107,0,1,98,0,0,8)
    ,@('Ruanda',105,0,7,98,0,0,0)
    ,@('Madagascar',93,5,11,82,1,0,0)
    ,@('El Salvador',93,15,9,79,2,1,5)
    ,@('Guatemala',80,3,17,60,3,0,3)
    ,@('Monaco',79,0,4,74,4,0,1)
    ,@('Liechtenstein',78,0,55,22,0,0,1)
    ,@('Guayana Francesa',77,0,34,43,1,0,0)
    ,@('Aruba',74,0,14,60,0,0,0)
    ,@('Togo',70,5,23,44,0,0,3)
    ,@('Barbados',63,0,6,54,4,0,3)
    ,@('Jamaica',63,0,10,50,0,0,3)
    ,@('Mali',56,0,12,39,0,0,5)
    ,@('Etiopia',55,3,4,49,1,0,2)
    ,@('Uganda',52,0,0,52,0,0,0)
    ,@('Polinesia Francesa',47,0,0,47,0,0,0)
    ,@('Congo',45,0,2,38,0,0,5)
    ,@('Islas Caimanes',45,0,6,38,0,0,1)
    ,@('Macao',45,1,10,35,1,0,0)
    ,@('San Martin (Parte Holandesa)',40,0,1,33,2,0,6)
    ,@('Puerto Rico',39,0,1,36,0,0,2)
    ,@('Zambia',39,0,7,31,0,0,1)
    ,@('Bermudas',39,0,17,20,0,0,2)
    ,@('Bahamas',36,3,5,25,1,0,6)
    ,@('Guinea-Bisau',33,0,0,33,0,0,0)
    ,@('Guyana',33,0,8,20,4,0,5)
    ,@('Guam',32,0,0,31,0,0,1)
    ,@('San Martin (Parte Francesa)',32,0,7,23,6,0,2)
    ,@('Eritrea',31,0,0,31,0,0,0)
    ,@('Gabon',30,0,1,28,0,0,1)
    ,@('Benin',26,0,5,20,0,0,1)
    ,@('Haiti',25,0,0,24,0,0,1)
    ,@('Tanzania',24,0,5,18,0,0,1)
    ,@('Birmania',22,0,0,19,0,2,3)
    ,@('Libia',21,1,2,18,0,0,1)
    ,@('Antigua y Barbuda',19,0,0,17,1,1,2)
    ,@('Siria',19,0,3,14,0,0,2)
    ,@('Maldivas',19,0,13,6,0,0,0)
    ,@('Nueva Caledonia',18,0,1,17,0,0,0)
    ,@('Islas Virgenes de los Estados Unidos',17,0,0,17,0,0,0)
    ,@('Angola',17,0,2,13,0,0,2)
    ,@('Namibia',16,0,3,13,0,0,0)
    ,@('Guinea Ecuatorial',16,0,3,13,0,0,0)
    ,@('Mongolia',16,1,4,12,0,0,0)
    ,@('Fiyi',15,0,0,15,0,0,0)
    ,@('Laos',15,1,0,15,0,0,0)
    ,@('Dominica',15,0,1,14,0,0,0)
    ,@('Santa Lucia',14,0,1,13,0,0,0)
    ,@('Sudan',14,0,2,10,0,0,2)
    ,@('Liberia',14,0,3,8,0,0,3)
    ,@('Curazao',13,0,7,5,0,0,1)
    ,@('Granada',12,0,0,12,2,0,0)
    ,@('San Cristobal y Nieves',11,0,0,11,0,0,0)
    ,@('Seychelles',11,0,0,11,0,0,0)
    ,@('Zimbabue',11,0,0,9,0,0,2)
    ,@('Groenlandia',11,0,10,1,0,0,0)
    ,@('Mozambique',10,0,1,9,0,0,0)
    ,@('Republica del Chad',10,0,2,8,0,0,0)
    ,@('Suazilandia',10,0,4,6,0,0,0)
    ,@('Surinam',10,0,3,6,0,0,1)
    ,@('Nepal',9,0,1,8,0,0,0)
    ,@('Montserrat',9,0,0,7,0,0,2)
    ,@('Republica de Africa Central',8,0,0,8,0,0,0)
    ,@('San Vicente y las Granadinas',8,0,1,7,0,0,0)
    ,@('Islas Turcas y Caicos',8,0,0,7,0,0,1)
    ,@('Somalia',8,0,1,7,0,0,0)
    ,@('Malaui',8,0,0,7,1,0,1)
    ,@('Belice',8,1,0,7,1,0,1)
    ,@('Sierra Leona',7,1,0,7,0,0,0)
    ,@('Santa Sede',7,0,0,7,0,0,0)
    ,@('Cabo Verde',7,0,1,5,0,0,1)
    ,@('Nicaragua',6,0,0,5,0,0,1)
    ,@('San Bartolome',6,0,1,5,0,0,0)
    ,@('Botsuana',6,0,0,5,0,0,1)
    ,@('Mauritania',6,0,2,3,0,0,1)
    ,@('Butan',5,0,2,3,0,0,0)
    ,@('Sahara Occidental',4,0,0,4,0,0,0)
    ,@('Santo Tome y Principe',4,0,0,4,0,0,0)
    ,@('Gambia',4,0,2,1,0,0,1)
    ,@('Anguila',3,0,0,3,0,0,0)
    ,@('Islas Virgenes Britanicas',3,0,0,3,0,0,0)
    ,@('Burundi',3,0,0,3,0,0,0)
    ,@('Papua Nueva Guinea',2,0,0,2,0,0,0)
    ,@('Islas Malvinas',2,0,0,2,0,0,0)
    ,@('Bonaire, San Eustaquio y Saba',2,0,0,2,0,0,0)
    ,@('Sudan del Sur',2,0,0,2,0,0,0)
    ,@('Timor Oriental',1,0,0,1,0,0,0)
    ,@('San Pedro y Miquelon',1,0,0,1,0,0,0)
)

$startRow = 4
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $ws.Cells.Item($row, 8).Value = $rec[7]
}
